# Update gh-pages to output generated at 456a3b4
# Applies numeric "want-to-go" count bumps and date-range text fixes
# across the four worksheets: 展览 (1), 演出 (2), 本地生活 (3), 全部类型 (4)

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item(1)   # 展览
$wsShow       = $wb.Worksheets.Item(2)   # 演出
$wsLocalLife  = $wb.Worksheets.Item(3)   # 本地生活
$wsAllTypes   = $wb.Worksheets.Item(4)   # 全部类型

# --- 展览 (sheet 1) ---
$wsExhibition.Range("F6").Value  = 547
$wsExhibition.Range("F10").Value = 381
$wsExhibition.Range("F12").Value = 657
$wsExhibition.Range("F13").Value = 752
$wsExhibition.Range("F16").Value = 883
$wsExhibition.Range("F20").Value = 305
$wsExhibition.Range("F23").Value = 99
$wsExhibition.Range("F24").Value = 6573
$wsExhibition.Range("F25").Value = 4910
$wsExhibition.Range("F29").Value = 161
$wsExhibition.Range("F33").Value = 192
$wsExhibition.Range("F34").Value = 246
$wsExhibition.Range("F38").Value = 244

# --- 演出 (sheet 2) ---
$wsShow.Range("E2").Value = "2024.08.02 19:30-10.27 16:55"

# --- 本地生活 (sheet 3) ---
$wsLocalLife.Range("F3").Value = 2450
$wsLocalLife.Range("F5").Value = 51

# --- 全部类型 (sheet 4) ---
$wsAllTypes.Range("E5").Value  = "2024.08.02 19:30-10.27 16:55"
$wsAllTypes.Range("F7").Value  = 51
$wsAllTypes.Range("F9").Value  = 547
$wsAllTypes.Range("F14").Value = 381
$wsAllTypes.Range("F16").Value = 657
$wsAllTypes.Range("F17").Value = 752
$wsAllTypes.Range("F20").Value = 883
$wsAllTypes.Range("F24").Value = 305
$wsAllTypes.Range("F26").Value = 99
$wsAllTypes.Range("F29").Value = 6573
$wsAllTypes.Range("F30").Value = 4910
$wsAllTypes.Range("F34").Value = 192
$wsAllTypes.Range("F35").Value = 246
$wsAllTypes.Range("F43").Value = 244
